$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L8").Value = 0.03512404007820401
$ws.Range("C9").Value = 0.04203954205791236
$ws.Range("A10").Value = 0.02981055223042308
$ws.Range("B10").Value = 0.03965627658747125
$ws.Range("C10").Value = 0.04163801216386761
$ws.Range("D10").Value = 0.04214682721648869
$ws.Range("E10").Value = 0.04175318980681017
$ws.Range("F10").Value = 0.04045105685498676
$ws.Range("G10").Value = 0.03909522527921545
$ws.Range("H10").Value = 0.03865774304801105
$ws.Range("I10").Value = 0.03818794853771629
$ws.Range("J10").Value = 0.0371882216317551
$ws.Range("K10").Value = 0.03577007278476298
$ws.Range("L10").Value = 0.03530461514376117
$ws.Range("M10").Value = 0.03408331784854801
$ws.Range("N10").Value = 0.03461888564544432
$ws.Range("O10").Value = 0.05948296948722619
$ws.Range("P10").Value = 0.07193080535559436
$ws.Range("Q10").Value = 0.07290062380123619
$ws.Range("R10").Value = 0.07228968475931551
$ws.Range("S10").Value = 0.07254554499864518
$ws.Range("T10").Value = 0.06874895844817834
$ws.Range("A11").Value = 0.02932228549105412
$ws.Range("B11").Value = 0.03984764296168757
$ws.Range("C11").Value = 0.04186578508891921
$ws.Range("D11").Value = 0.04189569415510519
$ws.Range("E11").Value = 0.04157940804560083
$ws.Range("F11").Value = 0.04033649575817468
$ws.Range("G11").Value = 0.038873063466944
$ws.Range("H11").Value = 0.03848278365276686
$ws.Range("I11").Value = 0.03828519165925346
$ws.Range("J11").Value = 0.03727591843567041
$ws.Range("K11").Value = 0.03567816130210945
$ws.Range("L11").Value = 0.03556955736915279
$ws.Range("M11").Value = 0.03388301946643414
$ws.Range("N11").Value = 0.03458231858363352
$ws.Range("O11").Value = 0.05923304362540343
$ws.Range("P11").Value = 0.07087891697282364
$ws.Range("Q11").Value = 0.07258383751452631
$ws.Range("R11").Value = 0.07191402695244994
$ws.Range("S11").Value = 0.07175975062622075
$ws.Range("T11").Value = 0.06850281333107071
$ws.Range("A12").Value = 0.02992784945695767
$ws.Range("B12").Value = 0.03950021223192253
$ws.Range("C12").Value = 0.04223306736861817
$ws.Range("D12").Value = 0.04198572691094003
$ws.Range("E12").Value = 0.04172338993320767
$ws.Range("F12").Value = 0.04061577231239526
$ws.Range("G12").Value = 0.03916751140102503
$ws.Range("H12").Value = 0.03902982157001228
$ws.Range("I12").Value = 0.03859057569916555
$ws.Range("J12").Value = 0.03733385581717758
$ws.Range("K12").Value = 0.03577160914090414
$ws.Range("L12").Value = 0.03567528164044301
$ws.Range("M12").Value = 0.03423204542047387
$ws.Range("N12").Value = 0.034700272720206
$ws.Range("O12").Value = 0.05973685582122366
$ws.Range("P12").Value = 0.07227757165249807
$ws.Range("Q12").Value = 0.07287464463656372
$ws.Range("R12").Value = 0.07207709317290208
$ws.Range("S12").Value = 0.0717755078077453
$ws.Range("T12").Value = 0.06872032540536609
